$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A4").Value = "Roshni"
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 10

$ws.Range("A5").Value = "Jyoti"
$ws.Range("B5").Value = 6
$ws.Range("C5").Value = 18

$ws.Range("C5").Select()
